$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")

# Insert a new row above the old row 5 ("venue.1") so the new
# "timezone" entry lands right after "location" in the tournament table.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "timezone"
$ws.Range("B5").Value = "Australia/Sydney"

# The sheet's table ("tournament") needs to grow to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I15"))

# Reflect the author's final on-screen state: Tournament tab active,
# with B5 (the newly entered timezone value) selected.
$ws.Activate()
$ws.Range("B5").Select() | Out-Null
